# Edit EBEWE Dates table to add reissued due dates.
#
# For the "0 or 1" and "2 or 3" LADBS-digit rows, the Initial Compliance
# Due Date and Initial Comparative Period cells get a second, bold line
# noting the reissued due date ("Sept 7, 2023*" / "Sept 1, 2018 -
# Sept 1, 2023*"). Those cells wrap text and the rows grow to fit two
# lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TwoLineCell {
    param(
        [string]$Addr,
        [string]$FirstLine,
        [string]$SecondLine
    )

    $top = $FirstLine + "`n"
    $full = $top + $SecondLine

    $rng = $ws.Range($Addr)
    $rng.Value = $full
    $rng.Characters($top.Length + 1, $SecondLine.Length).Font.Bold = $true
    $rng.WrapText = $true
}

# Column B (Initial Compliance Due Date) first, then column C (Initial
# Comparative Period) - matches the order the cells were authored in.
Set-TwoLineCell "B2" "Dec 1, 2021" "Sept 7, 2023*"
Set-TwoLineCell "B3" "Dec 1, 2022" "Sept 7, 2023*"
Set-TwoLineCell "C2" "Dec 1, 2016 - Dec 1, 2021" "Sept 1, 2018 - Sept 1, 2023*"
Set-TwoLineCell "C3" "Dec 1, 2017 - Dec 1, 2022" "Sept 1, 2018 - Sept 1, 2023*"

# Grow rows 2 & 3 so both lines are visible.
$ws.Rows.Item(2).RowHeight = 43.2
$ws.Rows.Item(3).RowHeight = 43.2

# Match the author's final selection.
$ws.Range("C4").Select() | Out-Null
